$d = $word.ActiveDocument

$replacements = @(
    @("56×24=", "44×83="),
    @("47×26=", "37×91="),
    @("18×40=", "55×71="),
    @("72×83=", "42×58="),
    @("88×41=", "88×77="),
    @("97×31=", "77×94="),
    @("46×24=", "11×88="),
    @("58×59=", "96×59="),
    @("82×17=", "33×68="),
    @("42×43=", "95×91="),
    @("20×34=", "90×59="),
    @("87×56=", "26×98="),
    @("95×83=", "37×87="),
    @("78×95=", "25×69="),
    @("85×14=", "31×30="),
    @("30×11=", "15×12="),
    @("72×96=", "19×99="),
    @("40×98=", "21×46="),
    @("97×34=", "47×84="),
    @("32×13=", "39×42="),
    @("12×94=", "20×53="),
    @("67×84=", "53×75="),
    @("65×25=", "93×28="),
    @("28×55=", "26×65="),
    @("88×27=", "76×45=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
